# Logboek Joey.xlsx — "add predictions + logboek joey"
#
# Week 2 of the logbook is reorganised: the existing Tuesday (row 19) and
# the trailing row (row 20) make way for a fuller Monday-Thursday writeup
# with a new Thursday entry (rows 24-26) describing work with Kaan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old trailing row (row 20) entirely -------------------------
# Its data (0.625 / 330 / "aan c# voor applicatie gewerkt" / "morgen kijken
# noël en Ik verder naar een oplossing") is recreated further down at row 22.
$ws.Range("C20:N20").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# --- New row 21: "woensdag" ------------------------------------------------
$ws.Range("B21").Value = "woensdag"
$ws.Range("G21").Value = "aan php voor website gewerkt"

# --- Row 17: Monday header, lowercase now ----------------------------------
$ws.Range("B17").Value = "maandag"

# --- Row 19: re-purposed as "dinsdag" with a new activity ------------------
$ws.Range("B19").Value = "dinsdag"

# --- New row 24: "Dondaerdag" (Thursday) ------------------------------------
$ws.Range("B24").Value = "Dondaerdag"

# --- Activity text for row 19 -----------------------------------------------
$ws.Range("G19").Value = "gekeken met noël naar wireframes"

# --- Activity text shared by rows 24 and 26 ---------------------------------
$ws.Range("G24").Value = "aan c# voor applicate gewerkt"

# --- Row 25: conversation with Kaan -----------------------------------------
$ws.Range("G25").Value = "kaan gevraagd wat hij deed"
$ws.Range("N25").Value = "we hebben hem een html taak gegeven"

# --- Fill in the remaining (numeric / already-existing-text) cells ---------

# Row 19 time moves from 13:30 to 13:00
$ws.Range("C19").NumberFormat = "h:mm"
$ws.Range("C19").Value = 0.54166666666666663
$ws.Range("E19").Value = 330

# Row 21 (woensdag)
$ws.Range("C21").NumberFormat = "h:mm"
$ws.Range("C21").Value = 0.5625
$ws.Range("E21").Value = 330

# Row 22 (re-created former row 20 content)
$ws.Range("C22").NumberFormat = "h:mm"
$ws.Range("C22").Value = 0.625
$ws.Range("E22").Value = 330
$ws.Range("G22").Value = "aan c# voor applicatie gewerkt"
$ws.Range("N22").Value = "morgen kijken noël en Ik verder naar een oplossing"

# Row 24 (Dondaerdag)
$ws.Range("C24").NumberFormat = "h:mm"
$ws.Range("C24").Value = 0.39583333333333331
$ws.Range("E24").Value = 330

# Row 25 (kaan gevraagd wat hij deed)
$ws.Range("C25").NumberFormat = "h:mm"
$ws.Range("C25").Value = 0.45833333333333331
$ws.Range("E25").Value = 330

# Row 26 (aan c# voor applicate gewerkt, again)
$ws.Range("C26").NumberFormat = "h:mm"
$ws.Range("C26").Value = 0.45902777777777781
$ws.Range("E26").Value = 330
$ws.Range("G26").Value = "aan c# voor applicate gewerkt"

# --- View state: scrolled down a bit, selection on K30 ---------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("K30").Select()
